$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ------------------------------------------------------------------
# 1. settings sheet: drop the "form_id" column (B) entirely.
#    Shift the cell comments one column to the left first (column
#    delete does not relocate comments on its own), then delete the
#    now-redundant trailing comment and finally remove the column.
# ------------------------------------------------------------------
$commentVersion    = $ws2.Range("C1").Comment.Text()
$commentPages      = $ws2.Range("D1").Comment.Text()
$commentNamespaces = $ws2.Range("E1").Comment.Text()

[void]$ws2.Range("B1").Comment.Text($commentVersion)
[void]$ws2.Range("C1").Comment.Text($commentPages)
[void]$ws2.Range("D1").Comment.Text($commentNamespaces)
$ws2.Range("E1").Comment.Delete()

$ws2.Range("B1").EntireColumn.Delete()

# ------------------------------------------------------------------
# 2. settings sheet: move the selection highlight from A11 to B11
# ------------------------------------------------------------------
$ws2.Activate()
[void]$ws2.Range("B11").Select()

# ------------------------------------------------------------------
# 3. survey sheet: consolidate the fragmented conditional formatting
#    sqref ranges back into single contiguous ranges.
# ------------------------------------------------------------------
$fcs = $ws1.Cells.FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fc = $fcs.Item($i)
    $addr = $fc.AppliesTo.Address()
    if ($addr -like "*`$A`$*" -and $addr -like "*`$F`$9993*") {
        $fc.ModifyAppliesToRange($ws1.Range("A2:F9993"))
    } elseif ($addr -like "*`$C`$*" -and $addr -like "*`$C`$9993*") {
        $fc.ModifyAppliesToRange($ws1.Range("C2:C9993"))
    }
}

# ------------------------------------------------------------------
# 4. survey sheet: move the bottomRight pane's selection to A15 and
#    re-activate this sheet so it remains the active tab.
# ------------------------------------------------------------------
$ws1.Activate()
[void]$ws1.Range("A15").Select()
